# Version 2 with catch trials and check screens
#
# The deck currently has 2 slides:
#   1) "请稍作休息"            (please take a short break)
#   2) "按 空格键 继续" / "请稍作休息"   (press space bar to continue)
#
# We insert a brand-new "check screen" slide in front of the break slide:
#   1) "请等待实验员操作"       (please wait for the experimenter) <- NEW
#   2) "请稍作休息"             (unchanged, shifted down)
#   3) "按 空格键 继续" / "请稍作休息" (unchanged, shifted down)
#
# The new slide is produced by duplicating slide 1 (so it inherits the same
# background, textbox formatting/size and shape identity) and then editing
# the duplicate's text + vertical position.

$p = $ppt.ActivePresentation

$original = $p.Slides.Item(1)

# Duplicate slide 1; the new slide lands immediately after it (index 2).
$newSlide = $original.Duplicate().Item(1)

# Move the duplicate in front, so it becomes the new slide 1.
$newSlide.MoveTo(1)

# Edit the duplicated slide's textbox: new wording + lower position.
$shape = $newSlide.Shapes.Item(1)
$shape.TextFrame.TextRange.Text = "请等待实验员操作"

# Shape.Top is in points; the OOXML offset is in EMU (1 pt = 12700 EMU).
$shape.Top = 2875002 / 12700
